# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# matching the target OOXML diff (currentAveragePrice* / Leve profit columns).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2468683.2
$ws.Range("J17").Value = 2468683.2
$ws.Range("L17").Value = 7406049.600000001
$ws.Range("N17").Value = -7406385.600000001

$ws.Range("H18").Value = 642.3570999999999
$ws.Range("I18").Value = 638.53845
$ws.Range("J18").Value = 692
$ws.Range("K18").Value = 638.53845
$ws.Range("L18").Value = 692
$ws.Range("M18").Value = -354.53845
$ws.Range("N18").Value = -1260

$ws.Range("H131").Value = 2715
$ws.Range("I131").Value = 2510
$ws.Range("J131").Value = 3125
$ws.Range("K131").Value = 7530
$ws.Range("L131").Value = 9375
$ws.Range("M131").Value = -2490
$ws.Range("N131").Value = -19455

$ws.Range("H138").Value = 2422.2073
$ws.Range("I138").Value = 1091.1143
$ws.Range("J138").Value = 3413.4468
$ws.Range("K138").Value = 3273.3429
$ws.Range("L138").Value = 10240.3404
$ws.Range("M138").Value = 1866.6571
$ws.Range("N138").Value = -20520.3404

$ws.Range("H141").Value = 1500.9362
$ws.Range("I141").Value = 1127.6875
$ws.Range("K141").Value = 3383.0625
$ws.Range("M141").Value = 1796.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6023.4346
$ws.Range("I32").Value = 4621.8716
$ws.Range("J32").Value = 11229.238
$ws.Range("K32").Value = 4621.8716
$ws.Range("L32").Value = 11229.238
$ws.Range("M32").Value = -4334.8716
$ws.Range("N32").Value = -11803.238

$ws.Range("H102").Value = 4116844
$ws.Range("I102").Value = 4631199.5
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 4631199.5
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -4629577.5
$ws.Range("N102").Value = -5244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1424.7894
$ws.Range("I94").Value = 512.3333
$ws.Range("K94").Value = 512.3333
$ws.Range("M94").Value = -61.33330000000001

$ws.Range("H99").Value = 45455716
$ws.Range("I99").Value = 55556470
$ws.Range("K99").Value = 55556470
$ws.Range("M99").Value = -55554972

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1361.4783
$ws.Range("I16").Value = 1262.7273
$ws.Range("K16").Value = 1262.7273
$ws.Range("M16").Value = -975.7273

$ws.Range("H31").Value = 269869.78
$ws.Range("I31").Value = 1675.2683
$ws.Range("K31").Value = 1675.2683
$ws.Range("M31").Value = -1380.2683

$ws.Range("H34").Value = 269869.78
$ws.Range("I34").Value = 1675.2683
$ws.Range("K34").Value = 1675.2683
$ws.Range("M34").Value = -1473.2683

$ws.Range("H41").Value = 7686.3335
$ws.Range("I41").Value = 1059
$ws.Range("J41").Value = 11000
$ws.Range("K41").Value = 1059
$ws.Range("L41").Value = 11000
$ws.Range("M41").Value = -631
$ws.Range("N41").Value = -11856

$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21250

$ws.Range("H51").Value = 21000
$ws.Range("J51").Value = 21000
$ws.Range("L51").Value = 21000
$ws.Range("N51").Value = -22472

$ws.Range("H58").Value = 1496.4166
$ws.Range("I58").Value = 886.2069
$ws.Range("J58").Value = 2427.7896
$ws.Range("K58").Value = 886.2069
$ws.Range("L58").Value = 2427.7896
$ws.Range("M58").Value = -683.2069
$ws.Range("N58").Value = -2833.7896

$ws.Range("H59").Value = 27675.8
$ws.Range("J59").Value = 27675.8
$ws.Range("L59").Value = 27675.8
$ws.Range("N59").Value = -29965.8

$ws.Range("H61").Value = 21000
$ws.Range("J61").Value = 21000
$ws.Range("L61").Value = 21000
$ws.Range("N61").Value = -21696

$ws.Range("H68").Value = 19600
$ws.Range("J68").Value = 29700
$ws.Range("L68").Value = 29700
$ws.Range("N68").Value = -31198

$ws.Range("H71").Value = 19600
$ws.Range("J71").Value = 29700
$ws.Range("L71").Value = 89100
$ws.Range("N71").Value = -96588

$ws.Range("H74").Value = 30250
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 30250
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 30250
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -31998

$ws.Range("H77").Value = 30250
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 30250
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 90750
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -99486

$ws.Range("H113").Value = 1361.4783
$ws.Range("I113").Value = 1262.7273
$ws.Range("K113").Value = 1262.7273
$ws.Range("M113").Value = 907.2727

$ws.Range("H132").Value = 1943.4755
$ws.Range("I132").Value = 1502.6809
$ws.Range("J132").Value = 3423.2856
$ws.Range("K132").Value = 4508.0427
$ws.Range("L132").Value = 10269.8568
$ws.Range("M132").Value = -1978.0427
$ws.Range("N132").Value = -15329.8568

$ws.Range("H134").Value = 1996.8182
$ws.Range("I134").Value = 2458.0286
$ws.Range("J134").Value = 1189.7
$ws.Range("K134").Value = 7374.085800000001
$ws.Range("L134").Value = 3569.1
$ws.Range("M134").Value = -4839.085800000001
$ws.Range("N134").Value = -8639.1

$ws.Range("H136").Value = 1496.4166
$ws.Range("I136").Value = 886.2069
$ws.Range("J136").Value = 2427.7896
$ws.Range("K136").Value = 2658.6207
$ws.Range("L136").Value = 7283.3688
$ws.Range("M136").Value = -108.6206999999999
$ws.Range("N136").Value = -12383.3688

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 35717260
$ws.Range("J75").Value = 35717260
$ws.Range("L75").Value = 107151780
$ws.Range("N75").Value = -107153776

$ws.Range("H78").Value = 35717260
$ws.Range("J78").Value = 35717260
$ws.Range("L78").Value = 321455340
$ws.Range("N78").Value = -321465324

$ws.Range("H113").Value = 200543.95
$ws.Range("I113").Value = 555.75
$ws.Range("J113").Value = 455074.4
$ws.Range("K113").Value = 1667.25
$ws.Range("L113").Value = 1365223.2
$ws.Range("M113").Value = 502.75
$ws.Range("N113").Value = -1369563.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2429.25
$ws.Range("I132").Value = 2039.4
$ws.Range("J132").Value = 2819.1
$ws.Range("K132").Value = 6118.200000000001
$ws.Range("L132").Value = 8457.299999999999
$ws.Range("M132").Value = -3588.200000000001
$ws.Range("N132").Value = -13517.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11373068
$ws.Range("I132").Value = 14065695
$ws.Range("J132").Value = 4200.3335
$ws.Range("K132").Value = 42197085
$ws.Range("L132").Value = 12601.0005
$ws.Range("M132").Value = -42194555
$ws.Range("N132").Value = -17661.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 28950
$ws.Range("J123").Value = 28950
$ws.Range("L123").Value = 28950
$ws.Range("N123").Value = -38750

$ws.Range("H136").Value = 8477237
$ws.Range("I136").Value = 2934.4722
$ws.Range("J136").Value = 21741362
$ws.Range("K136").Value = 8803.4166
$ws.Range("L136").Value = 65224086
$ws.Range("M136").Value = -6253.4166
$ws.Range("N136").Value = -65229186
